$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Cells.Item(2, 2).Value = 1.426125737049745
$ws.Cells.Item(2, 3).Value = 0.2716980152002861
$ws.Cells.Item(2, 4).Value = 0.02815277586733345
$ws.Cells.Item(2, 5).Value = 0.09735005301325117
$ws.Cells.Item(2, 6).Value = 0.7404525498735524
$ws.Cells.Item(2, 8).Value = 0.07973214163530429
$ws.Cells.Item(2, 12).Value = 0.2186133324869957
$ws.Cells.Item(2, 15).Value = 2.545324754224538

# Row 3
$ws.Cells.Item(3, 2).Value = 1.287646044420967
$ws.Cells.Item(3, 3).Value = 0.2566188330006582
$ws.Cells.Item(3, 4).Value = 0.02670168647993165
$ws.Cells.Item(3, 5).Value = 0.09856702563679853
$ws.Cells.Item(3, 6).Value = 0.7407862203640079
$ws.Cells.Item(3, 8).Value = 0.07973214163530429
$ws.Cells.Item(3, 12).Value = 0.2082890078713007
$ws.Cells.Item(3, 15).Value = 2.562567931480714

# Row 4
$ws.Cells.Item(4, 2).Value = 1.202640181683364
$ws.Cells.Item(4, 3).Value = 0.2473352152662471
$ws.Cells.Item(4, 4).Value = 0.02580839362496334
$ws.Cells.Item(4, 5).Value = 0.09937666743672935
$ws.Cells.Item(4, 6).Value = 0.7416464198107562
$ws.Cells.Item(4, 8).Value = 0.07973214163530429
$ws.Cells.Item(4, 12).Value = 0.202049473344843
$ws.Cells.Item(4, 15).Value = 2.575409500003786

# Row 5
$ws.Cells.Item(5, 2).Value = 1.168006751467374
$ws.Cells.Item(5, 3).Value = 0.2435460889232388
$ws.Cells.Item(5, 4).Value = 0.02544381431792431
$ws.Cells.Item(5, 5).Value = 0.0997222891822922
$ws.Cells.Item(5, 6).Value = 0.742161351776808
$ws.Cells.Item(5, 8).Value = 0.07973214163530429
$ws.Cells.Item(5, 12).Value = 0.1995319262775155
$ws.Cells.Item(5, 15).Value = 2.581207977441551

# Row 6
$ws.Cells.Item(6, 2).Value = 1.162256382370003
$ws.Cells.Item(6, 3).Value = 0.2429165551416475
$ws.Cells.Item(6, 4).Value = 0.02538324349009713
$ws.Cells.Item(6, 5).Value = 0.09978062655603814
$ws.Cells.Item(6, 6).Value = 0.7422567732847298
$ws.Cells.Item(6, 8).Value = 0.07973214163530429
$ws.Cells.Item(6, 12).Value = 0.1991154080678683
$ws.Cells.Item(6, 15).Value = 2.58220491995786

# Row 7
$ws.Cells.Item(7, 2).Value = 1.202173071766651
$ws.Cells.Item(7, 3).Value = 0.2472841375883092
$ws.Cells.Item(7, 4).Value = 0.02580347899287716
$ws.Cells.Item(7, 5).Value = 0.09938126511155332
$ws.Cells.Item(7, 6).Value = 0.7416526992634402
$ws.Cells.Item(7, 8).Value = 0.07973214163530429
$ws.Cells.Item(7, 12).Value = 0.2020154190604018
$ws.Cells.Item(7, 15).Value = 2.57548541274511

# Row 8
$ws.Cells.Item(8, 2).Value = 1.378374569437199
$ws.Cells.Item(8, 3).Value = 0.2665040599691224
$ws.Cells.Item(8, 4).Value = 0.02765293741384767
$ws.Cells.Item(8, 5).Value = 0.09775670556903471
$ws.Cells.Item(8, 6).Value = 0.7404313483539795
$ws.Cells.Item(8, 8).Value = 0.07973214163530429
$ws.Cells.Item(8, 12).Value = 0.2150328553357497
$ws.Cells.Item(8, 15).Value = 2.55080149997832

# Row 9
$ws.Cells.Item(9, 2).Value = 1.724012180240891
$ws.Cells.Item(9, 3).Value = 0.3039854258960872
$ws.Cells.Item(9, 4).Value = 0.03126030033718763
$ws.Cells.Item(9, 5).Value = 0.09506657982693589
$ws.Cells.Item(9, 6).Value = 0.7432548820247291
$ws.Cells.Item(9, 8).Value = 0.07973214163530429
$ws.Cells.Item(9, 12).Value = 0.241349659559603
$ws.Cells.Item(9, 15).Value = 2.520347592924026

# Row 10
$ws.Cells.Item(10, 2).Value = 1.977961545731148
$ws.Cells.Item(10, 3).Value = 0.3313838256891586
$ws.Cells.Item(10, 4).Value = 0.03389769180528646
$ws.Cells.Item(10, 5).Value = 0.09339282582797992
$ws.Cells.Item(10, 6).Value = 0.748538571820788
$ws.Cells.Item(10, 8).Value = 0.07973214163530429
$ws.Cells.Item(10, 12).Value = 0.2611670464069249
$ws.Cells.Item(10, 15).Value = 2.509009849006873

# Row 11
$ws.Cells.Item(11, 2).Value = 2.09348190978568
$ws.Cells.Item(11, 3).Value = 0.3438155814336881
$ws.Cells.Item(11, 4).Value = 0.03509449064593184
$ws.Cells.Item(11, 5).Value = 0.09269723182196898
$ws.Cells.Item(11, 6).Value = 0.7516451862771873
$ws.Cells.Item(11, 8).Value = 0.07973214163530429
$ws.Cells.Item(11, 12).Value = 0.2702876302102482
$ws.Cells.Item(11, 15).Value = 2.506268740626041

# Row 12
$ws.Cells.Item(12, 2).Value = 2.137224739684825
$ws.Cells.Item(12, 3).Value = 0.3485183309733202
$ws.Cells.Item(12, 4).Value = 0.03554723832923656
$ws.Cells.Item(12, 5).Value = 0.09244330170096049
$ws.Cells.Item(12, 6).Value = 0.7529231317852805
$ws.Cells.Item(12, 8).Value = 0.07973214163530429
$ws.Cells.Item(12, 12).Value = 0.2737565190599156
$ws.Cells.Item(12, 15).Value = 2.505579824533129

# Row 13
$ws.Cells.Item(13, 2).Value = 2.127804067047634
$ws.Cells.Item(13, 3).Value = 0.3475057317340031
$ws.Cells.Item(13, 4).Value = 0.03544975164792419
$ws.Cells.Item(13, 5).Value = 0.09249756842256218
$ws.Cells.Item(13, 6).Value = 0.7526433804024464
$ws.Cells.Item(13, 8).Value = 0.07973214163530429
$ws.Cells.Item(13, 12).Value = 0.2730087601340756
$ws.Cells.Item(13, 15).Value = 2.505712644077278

# Row 14
$ws.Cells.Item(14, 2).Value = 2.097080713932883
$ws.Cells.Item(14, 3).Value = 0.3442025790885452
$ws.Cells.Item(14, 4).Value = 0.03513174772825067
$ws.Cells.Item(14, 5).Value = 0.09267615084283776
$ws.Cells.Item(14, 6).Value = 0.7517482861240268
$ws.Cells.Item(14, 8).Value = 0.07973214163530429
$ws.Cells.Item(14, 12).Value = 0.2705727150446222
$ws.Cells.Item(14, 15).Value = 2.506205057842976

# Row 15
$ws.Cells.Item(15, 2).Value = 2.078261428250983
$ws.Cells.Item(15, 3).Value = 0.3421786575931662
$ws.Cells.Item(15, 4).Value = 0.03493690119242387
$ws.Cells.Item(15, 5).Value = 0.09278677221847786
$ws.Cells.Item(15, 6).Value = 0.7512132513722634
$ws.Cells.Item(15, 8).Value = 0.07973214163530429
$ws.Cells.Item(15, 12).Value = 0.2690825347070529
$ws.Cells.Item(15, 15).Value = 2.506552183530147

# Row 16
$ws.Cells.Item(16, 2).Value = 1.970411808117262
$ws.Cells.Item(16, 3).Value = 0.3305707135031355
$ws.Cells.Item(16, 4).Value = 0.03381941615469231
$ws.Cells.Item(16, 5).Value = 0.09343960985674116
$ws.Cells.Item(16, 6).Value = 0.7483497320056358
$ws.Cells.Item(16, 8).Value = 0.07973214163530429
$ws.Cells.Item(16, 12).Value = 0.2605731139490786
$ws.Cells.Item(16, 15).Value = 2.509237749174559

# Row 17
$ws.Cells.Item(17, 2).Value = 1.904247587643624
$ws.Cells.Item(17, 3).Value = 0.3234412288711042
$ws.Cells.Item(17, 4).Value = 0.03313309620906324
$ws.Cells.Item(17, 5).Value = 0.09385696850501368
$ws.Cells.Item(17, 6).Value = 0.7467734355434317
$ws.Cells.Item(17, 8).Value = 0.07973214163530429
$ws.Cells.Item(17, 12).Value = 0.2553798542435288
$ws.Cells.Item(17, 15).Value = 2.511505312961702

# Row 18
$ws.Cells.Item(18, 2).Value = 1.866191555075147
$ws.Cells.Item(18, 3).Value = 0.3193375514979664
$ws.Cells.Item(18, 4).Value = 0.03273806580221361
$ws.Cells.Item(18, 5).Value = 0.09410321523864162
$ws.Cells.Item(18, 6).Value = 0.7459329384071367
$ws.Cells.Item(18, 8).Value = 0.07973214163530429
$ws.Cells.Item(18, 12).Value = 0.2524027673105707
$ws.Cells.Item(18, 15).Value = 2.513036911941953

# Row 19
$ws.Cells.Item(19, 2).Value = 1.85330647257706
$ws.Cells.Item(19, 3).Value = 0.3179476130713965
$ws.Cells.Item(19, 4).Value = 0.03260426852084919
$ws.Cells.Item(19, 5).Value = 0.09418765340573465
$ws.Cells.Item(19, 6).Value = 0.7456597078140987
$ws.Cells.Item(19, 8).Value = 0.07973214163530429
$ws.Cells.Item(19, 12).Value = 0.2513964855610311
$ws.Cells.Item(19, 15).Value = 2.513594485067813

# Row 20
$ws.Cells.Item(20, 2).Value = 1.911290908205387
$ws.Cells.Item(20, 3).Value = 0.3242004859288272
$ws.Cells.Item(20, 4).Value = 0.03320618504430684
$ws.Cells.Item(20, 5).Value = 0.09381189890598129
$ws.Cells.Item(20, 6).Value = 0.7469343858267337
$ws.Cells.Item(20, 8).Value = 0.07973214163530429
$ws.Cells.Item(20, 12).Value = 0.2559316574052843
$ws.Cells.Item(20, 15).Value = 2.511240384834565

# Row 21
$ws.Cells.Item(21, 2).Value = 2.106104982639408
$ws.Cells.Item(21, 3).Value = 0.345172930058169
$ws.Cells.Item(21, 4).Value = 0.03522516573077183
$ws.Cells.Item(21, 5).Value = 0.09262343956701002
$ws.Cells.Item(21, 6).Value = 0.7520084377834166
$ws.Cells.Item(21, 8).Value = 0.07973214163530429
$ws.Cells.Item(21, 12).Value = 0.2712878305881503
$ws.Cells.Item(21, 15).Value = 2.506050937310846

# Row 22
$ws.Cells.Item(22, 2).Value = 2.233413424287846
$ws.Cells.Item(22, 3).Value = 0.3588510172432109
$ws.Cells.Item(22, 4).Value = 0.03654202706910326
$ws.Cells.Item(22, 5).Value = 0.09190195369454557
$ws.Cells.Item(22, 6).Value = 0.7559166693299488
$ws.Cells.Item(22, 8).Value = 0.07973214163530429
$ws.Cells.Item(22, 12).Value = 0.281412092657277
$ws.Cells.Item(22, 15).Value = 2.504694788127921

# Row 23
$ws.Cells.Item(23, 2).Value = 2.165468388404918
$ws.Cells.Item(23, 3).Value = 0.3515534802873503
$ws.Cells.Item(23, 4).Value = 0.03583944580750398
$ws.Cells.Item(23, 5).Value = 0.09228196525754839
$ws.Cells.Item(23, 6).Value = 0.7537764548151671
$ws.Cells.Item(23, 8).Value = 0.07973214163530429
$ws.Cells.Item(23, 12).Value = 0.2760005394595026
$ws.Cells.Item(23, 15).Value = 2.505231820633924

# Row 24
$ws.Cells.Item(24, 2).Value = 1.908106676857585
$ws.Cells.Item(24, 3).Value = 0.3238572408721154
$ws.Cells.Item(24, 4).Value = 0.03317314299945195
$ws.Cells.Item(24, 5).Value = 0.09383225524345384
$ws.Cells.Item(24, 6).Value = 0.7468614154691693
$ws.Cells.Item(24, 8).Value = 0.07973214163530429
$ws.Cells.Item(24, 12).Value = 0.2556821604459003
$ws.Cells.Item(24, 15).Value = 2.511359448910099

# Row 25
$ws.Cells.Item(25, 2).Value = 1.630502460955199
$ws.Cells.Item(25, 3).Value = 0.2938693232786136
$ws.Cells.Item(25, 4).Value = 0.03028660745812317
$ws.Cells.Item(25, 5).Value = 0.09574121442087424
$ws.Cells.Item(25, 6).Value = 0.7419294698460845
$ws.Cells.Item(25, 8).Value = 0.07973214163530429
$ws.Cells.Item(25, 12).Value = 0.2341456260844126
$ws.Cells.Item(25, 15).Value = 2.526654732982195
